$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format specific Price cells as Text so numeric-looking strings
# (e.g. "1.005") are preserved verbatim instead of being parsed as numbers,
# matching the source data which stores prices as text.
$textCells = @("D5", "D6", "D7", "D8", "D9", "D10", "D11", "D13", "D14", "D15", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D39", "D40", "D41", "D42", "D43", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated coin price / 1h volume-change figures
$ws.Range("D2").Value = '25.986.56'
$ws.Range("E2").Value = '  -0.37%  '
$ws.Range("D3").Value = '1.642.99'
$ws.Range("E3").Value = '  -1.39%  '
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").Value = '205.04'
$ws.Range("E5").Value = '  -1.82%  '
$ws.Range("D6").Value = '0.5165'
$ws.Range("E6").Value = '  +0.14%  '
$ws.Range("D7").Value = '1.005'
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("D8").Value = '0.2561'
$ws.Range("E8").Value = '  -2.77%  '
$ws.Range("D9").Value = '0.06221'
$ws.Range("E9").Value = '  +0.19%  '
$ws.Range("D10").Value = '20.58'
$ws.Range("E10").Value = '  -1.90%  '
$ws.Range("D11").Value = '0.07546'
$ws.Range("E11").Value = '  +0.67%  '
$ws.Range("D12").Value = '1.648.22'
$ws.Range("E12").Value = '  -0.94%  '
$ws.Range("D13").Value = '4.362'
$ws.Range("E13").Value = '  -1.23%  '
$ws.Range("D14").Value = '0.5347'
$ws.Range("E14").Value = '  -4.02%  '
$ws.Range("D15").Value = '65.85'
$ws.Range("E15").Value = '  +0.70%  '
$ws.Range("D16").Value = '0.0₅7868'
$ws.Range("E16").Value = '  -0.77%  '
$ws.Range("D17").Value = '25.990.52'
$ws.Range("E17").Value = '  -0.42%  '
$ws.Range("D18").Value = '1.004'
$ws.Range("E18").Value = '  +0.10%  '
$ws.Range("D19").Value = '4.645'
$ws.Range("E19").Value = '  -2.73%  '
$ws.Range("D20").Value = '186.11'
$ws.Range("E20").Value = '  +0.02%  '
$ws.Range("D21").Value = '9.915'
$ws.Range("E21").Value = '  -4.32%  '
$ws.Range("D22").Value = '6.101'
$ws.Range("E22").Value = '  -0.45%  '
$ws.Range("D23").Value = '1.006'
$ws.Range("E23").Value = '  +0.25%  '
$ws.Range("D24").Value = '147.34'
$ws.Range("E24").Value = '  +0.51%  '
$ws.Range("D25").Value = '0.1206'
$ws.Range("E25").Value = '  -2.75%  '
$ws.Range("D26").Value = '7.292'
$ws.Range("E26").Value = '  -3.12%  '
$ws.Range("D27").Value = '15.57'
$ws.Range("E27").Value = '  -1.06%  '
$ws.Range("D28").Value = '1.386'
$ws.Range("E28").Value = '  +3.29%  '
$ws.Range("D29").Value = '0.05952'
$ws.Range("E29").Value = '  -5.74%  '
$ws.Range("D30").Value = '1.238'
$ws.Range("E30").Value = '  -2.37%  '
$ws.Range("D31").Value = '3.430'
$ws.Range("E31").Value = '  -1.14%  '
$ws.Range("D32").Value = '3.372'
$ws.Range("E32").Value = '  -1.41%  '
$ws.Range("D33").Value = '1.611'
$ws.Range("E33").Value = '  +0.31%  '
$ws.Range("D34").Value = '0.9682'
$ws.Range("E34").Value = '  -2.50%  '
$ws.Range("D35").Value = '2.388'
$ws.Range("E35").Value = '  -0.69%  '
$ws.Range("D36").Value = '2.720'
$ws.Range("E36").Value = '  +0.64%  '
$ws.Range("D37").Value = '0.5835'
$ws.Range("E37").Value = '  -2.79%  '
$ws.Range("D38").Value = '1.084.19'
$ws.Range("E38").Value = '  +0.42%  '
$ws.Range("D39").Value = '0.01581'
$ws.Range("E39").Value = '  -1.75%  '
$ws.Range("D40").Value = '5.853'
$ws.Range("E40").Value = '  -3.67%  '
$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D41").Value = '1.002'
$ws.Range("E41").Value = '  -0.12%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = '0.8401'
$ws.Range("E42").Value = '  -2.10%  '
$ws.Range("D43").Value = '100.11'
$ws.Range("E43").Value = '  +1.06%  '
$ws.Range("D44").Value = '1.805.28'
$ws.Range("E44").Value = '  -0.43%  '
$ws.Range("D45").Value = '0.0₈110'
$ws.Range("E45").Value = '  -0.26%  '
$ws.Range("B46").Value = 'Frax'
$ws.Range("C46").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D46").Value = '0.9919'
$ws.Range("E46").Value = '  -1.32%  '
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = '54.47'
$ws.Range("E47").Value = '  -2.68%  '
$ws.Range("D48").Value = '7.966'
$ws.Range("E48").Value = '  +0.78%  '
$ws.Range("D49").Value = '0.05218'
$ws.Range("E49").Value = '  -0.47%  '
$ws.Range("D50").Value = '0.4238'
$ws.Range("E50").Value = '  -0.42%  '
$ws.Range("D51").Value = '5.803'
$ws.Range("E51").Value = '  -1.53%  '
